$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 12:22"

# 2) Refresh Rusia's (row 13) figures
$ws.Cells.Item(13, 2).Value = 74588
$ws.Cells.Item(13, 3).Value = 5966
$ws.Cells.Item(13, 4).Value = 6250
$ws.Cells.Item(13, 5).Value = 67657
$ws.Cells.Item(13, 6).Value = 2300
$ws.Cells.Item(13, 7).Value = 66
$ws.Cells.Item(13, 8).Value = 681

# 3) Refresh Rumania's (row 35) figures
$ws.Cells.Item(35, 2).Value = 10635
$ws.Cells.Item(35, 3).Value = 218
$ws.Cells.Item(35, 4).Value = 2890
$ws.Cells.Item(35, 5).Value = 7170
$ws.Cells.Item(35, 6).Value = 235
$ws.Cells.Item(35, 7).Value = 8
$ws.Cells.Item(35, 8).Value = 575

# 4) Refresh Marruecos' (row 55) figures
$ws.Cells.Item(55, 2).Value = 3889
$ws.Cells.Item(55, 3).Value = 131
$ws.Cells.Item(55, 4).Value = 498
$ws.Cells.Item(55, 5).Value = 3232
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 159

# 5) Refresh Estonia's (row 72) figures
$ws.Cells.Item(72, 2).Value = 1635
$ws.Cells.Item(72, 3).Value = 30
$ws.Cells.Item(72, 4).Value = 228
$ws.Cells.Item(72, 5).Value = 1361
$ws.Cells.Item(72, 6).Value = 6
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 46

# 6) Albania's update pushes it up the leaderboard (past Crucero, Libano,
#    Costa Rica and Niger), so rows 95-99 are re-sorted descending by
#    "Casos totales" (column B). Rewrite rows 95-99 with their new
#    country/figures in the new sorted order.
$ws.Cells.Item(95, 1).Value = "Albania"
$ws.Cells.Item(95, 2).Value = 712
$ws.Cells.Item(95, 3).Value = 34
$ws.Cells.Item(95, 4).Value = 403
$ws.Cells.Item(95, 5).Value = 282
$ws.Cells.Item(95, 6).Value = 4
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 27

$ws.Cells.Item(96, 1).Value = "Crucero"
$ws.Cells.Item(96, 2).Value = 712
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 645
$ws.Cells.Item(96, 5).Value = 54
$ws.Cells.Item(96, 6).Value = 4
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 13

$ws.Cells.Item(97, 1).Value = "Libano"
$ws.Cells.Item(97, 2).Value = 704
$ws.Cells.Item(97, 3).Value = 8
$ws.Cells.Item(97, 4).Value = 143
$ws.Cells.Item(97, 5).Value = 537
$ws.Cells.Item(97, 6).Value = 44
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = 24

$ws.Cells.Item(98, 1).Value = "Costa Rica"
$ws.Cells.Item(98, 2).Value = 687
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 216
$ws.Cells.Item(98, 5).Value = 465
$ws.Cells.Item(98, 6).Value = 7
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 6

$ws.Cells.Item(99, 1).Value = "Niger"
$ws.Cells.Item(99, 2).Value = 681
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 289
$ws.Cells.Item(99, 5).Value = 368
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 24

# Row 100 (Kirguistan) is unchanged.

# 7) Bonaire, San Eustaquio y Saba overtakes Sudan del Sur (tie on column B,
#    order swapped) - swap the two country labels in rows 210/211; the
#    figures are identical for both so no numeric change is required.
$ws.Cells.Item(210, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(211, 1).Value = "Sudan del Sur"
